# Apply "commit results and report" changes to the Branch_and_cut_sol workbook.
#
# The workbook records one run per input instance: elapsed time (B), the best
# objective value found (C, stored as text), the clusters of the best
# solution found (D, stored as text), the best bound (E) and the relative
# gap (F). This script refreshes the numbers/text with the newly reported
# results.

function Set-TextValue($rng, $text) {
    # Plain string assignment lets Excel auto-detect numeric-looking text
    # (e.g. "42.74...") and silently store it as a Number cell, which would
    # change the cell's type/shared-string usage versus the source data.
    # Forcing a text-concatenation formula guarantees a Text result, and
    # copy / paste-special (values only) bakes it back down to a static,
    # non-formula shared-string cell - matching the original cell shape.
    $rng.Formula = '=""&"' + $text + '"'
    $rng.Copy() | Out-Null
    $rng.PasteSpecial(-4163) | Out-Null
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 14_burma_6.tsp row (row 7) ---
Set-TextValue $ws.Range("C7") "42.740623542601746"
Set-TextValue $ws.Range("D7") "[7][6, 12, 14][2, 13][3, 4, 5][9, 10][1, 8, 11]"
$ws.Range("B7").Value = 30.019999980926514
$ws.Range("F7").Value = 0.9999999999976603

# --- 22_ulysses_3.tsp row (row 9) ---
Set-TextValue $ws.Range("C9") "414.58990456431127"
Set-TextValue $ws.Range("D9") "[5, 6, 7, 13, 14, 15, 18][1, 2, 3, 4, 8, 16, 17, 20, 22][9, 10, 11, 12, 19, 21]"
$ws.Range("B9").Value = 30.017000198364258
$ws.Range("E9").Value = 173.61775574493146
$ws.Range("F9").Value = 0.5812301413189429

# --- Remaining timing/metric updates (rows 3-12) ---
$ws.Range("B3").Value = 2.509000062942505
$ws.Range("B4").Value = 20.182999849319458
$ws.Range("B5").Value = 10.39300012588501
$ws.Range("B6").Value = 3.1100001335144043
$ws.Range("B8").Value = 30.013999938964844
$ws.Range("B10").Value = 30.027000188827515
$ws.Range("E10").Value = 177.98397630890767
$ws.Range("F10").Value = 0.9487909546152236
$ws.Range("B11").Value = 30.06599998474121
$ws.Range("B12").Value = 30.051000118255615
